# The presentation currently carries its custom "Integral" theme in
# ppt/theme/theme1.xml (the theme actually used by the slide master /
# every slide) and a default "Office Theme" in ppt/theme/theme2.xml
# (only wired to the notes master). The edit swaps the two themes so the
# deck's slide master now carries the default Office color palette.
#
# PowerPoint's automation model only exposes theme colors through
# Master.Theme.ThemeColorScheme (there is no supported way to rewrite a
# theme's <a:theme name="..."> / <a:clrScheme name="..."> labels, or to
# reach the notes-master's theme part directly) so we reproduce the
# substantive part of the change -- the 12 scheme colors -- on the
# slide master's theme, driving it from the current "Integral" palette
# to the standard "Office" palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB -> long packed as 0xBBGGRR, the COM ColorFormat convention.
$colors.Item(1).RGB  = 0x000000     # dk1      -> 000000
$colors.Item(2).RGB  = 0xFFFFFF     # lt1      -> FFFFFF
$colors.Item(3).RGB  = 0x6A5444     # dk2      -> 44546A
$colors.Item(4).RGB  = 0xE6E6E7     # lt2      -> E7E6E6
$colors.Item(5).RGB  = 0xD59B5B     # accent1  -> 5B9BD5
$colors.Item(6).RGB  = 0x317DED     # accent2  -> ED7D31
$colors.Item(7).RGB  = 0xA5A5A5     # accent3  -> A5A5A5
$colors.Item(8).RGB  = 0x00C0FF     # accent4  -> FFC000
$colors.Item(9).RGB  = 0xC47244     # accent5  -> 4472C4
$colors.Item(10).RGB = 0x47AD70     # accent6  -> 70AD47
$colors.Item(11).RGB = 0xC16305     # hlink    -> 0563C1
$colors.Item(12).RGB = 0x724F95     # folHlink -> 954F72
